$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "K2" = 3
    "L2" = 1
    "M2" = 7.107333666666666
    "N2" = 21.322001
    "O2" = 0.7373665550576455
    "P2" = 0.7373665550576454
    "Q2" = 0.3998988669774444
    "R2" = 3.599089802797
    "S2" = 0.2011998777819247
    "T2" = 0.2011998777819247

    "O3" = 0.1688878844614928
    "P3" = 0.1688878844614928
    "S3" = 0.04608321530102989
    "T3" = 0.04608321530102987

    "M4" = 0.8135026666666666
    "N4" = 2.440508
    "O4" = 0.08439869112428164
    "P4" = 0.08439869112428162
    "Q4" = 0.04577226987511111
    "R4" = 0.411950428876
    "S4" = 0.02302926030843961
    "T4" = 0.0230292603084396

    "K5" = 2
    "L5" = 0.6666666666666666
    "M5" = 0.09009266666666667
    "N5" = 0.270278
    "O5" = 0.009346869356580103
    "P5" = 0.009346869356580103
    "Q5" = 0.005069123951777778
    "R5" = 0.045622115566
    "S5" = 0.00255041262624193
    "T5" = 0.00255041262624193

    "K6" = 3
    "L6" = 1
    "M6" = 7.107333666666666
    "N6" = 21.322001
    "O6" = 0.7373665550576455
    "P6" = 0.7373665550576454
    "Q6" = 1.065668871757556
    "R6" = 9.591019845818002
    "S6" = 0.5361666772757208
    "T6" = 0.5361666772757206

    "O7" = 0.1688878844614928
    "P7" = 0.1688878844614928
    "S7" = 0.1228046691604629
    "T7" = 0.1228046691604629

    "M8" = 0.8135026666666666
    "N8" = 2.440508
    "O8" = 0.08439869112428164
    "P8" = 0.08439869112428162
    "Q8" = 0.1219760475048889
    "R8" = 1.097784427544
    "S8" = 0.06136943081584204
    "T8" = 0.06136943081584202

    "K9" = 2
    "L9" = 0.6666666666666666
    "M9" = 0.09009266666666667
    "N9" = 0.270278
    "O9" = 0.009346869356580103
    "P9" = 0.009346869356580103
    "Q9" = 0.01350843437822222
    "R9" = 0.121575909404
    "S9" = 0.006796456730338173
    "T9" = 0.006796456730338173
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
